$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the workbook's default ("Normal") style font in place, matching
# the target: size 11 -> 10, name Calibri -> Bahnschrift.
$normalFont = $wb.Styles.Item("Normal").Font
$normalFont.Size = 10
$normalFont.Name = "Bahnschrift"

# Add the two new shared-string values next to the existing "jc !" cell.
$ws.Range("B1").Value = "Simple Text"
$ws.Range("C1").Value = "This is PhpSpreadSheet"
